$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final-draft copy edits to the "Method" column (G) — typo/formatting fixes
# and a couple of shortened summaries, per the author's "final draft for
# submission" pass.
$ws.Range("G2").Value = 'TimeGpt (with and without predictors), LSTM (with and without predictors), Exponential Smoothing, ARIMA, Dynamic Regression (ARIMAX), Multiple Regression (with and without predictors)'
$ws.Range("G3").Value = 'ARIMA, HOLT WINTER, ETS,  Theta, ARHOW'
$ws.Range("G8").Value = 'Diffusion models,  ARIMA, ES, LR'
$ws.Range("G11").Value = 'MA, LR, Clustering and RNN'
$ws.Range("G13").Value = 'SES, MA, Naïve, Holt  '
$ws.Range("G18").Value = 'Naïve, Seasonal Naïve, ARIMA, DLM, LST, Recurrent Neural Network, SES, Holts linear and damped trends, HW, ES-SSM, Theta'
$ws.Range("G22").Value = 'Brown’s method, Holt’s method, and Holt-Winters method, moving average'
$ws.Range("G23").Value = 'Naïve, ARIMA, XGBoost , LSTM, Prophet '

# Restore the last-used selection reflected by the saved view state.
$ws.Range("G26").Select()
